$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Temperature" (C) and "Precipitation" (D) columns to remain
# text cells (they hold numeric-looking / percent-looking strings such as
# "63" and "1%") instead of being auto-coerced into numbers by the Value
# setter. Restoring NumberFormat to "" afterwards keeps the cell style
# identical to its original ("General") appearance.
$ws.Range("C2:D6").NumberFormat = "@"

# Row 2 - Kansas City
$ws.Range("A2").Value = "12/29/2022 10:46:49"
$ws.Range("B2").Value = "Kansas City"
$ws.Range("C2").Value = "63"
$ws.Range("D2").Value = "1%"
$ws.Range("E2").Value = "19 mph"
$ws.Range("F2").Value = "19 mph"

# Row 3 - New York
$ws.Range("A3").Value = "12/29/2022 10:46:52"
$ws.Range("B3").Value = "New York"
$ws.Range("C3").Value = "44"
$ws.Range("D3").Value = "1%"
$ws.Range("E3").Value = "6 mph"
$ws.Range("F3").Value = "19 mph"

# Row 4 - Sacramento
$ws.Range("A4").Value = "12/29/2022 10:46:55"
$ws.Range("B4").Value = "Sacramento"
$ws.Range("C4").Value = "40"
$ws.Range("D4").Value = "85%"
$ws.Range("E4").Value = "7 mph"
$ws.Range("F4").Value = "19 mph"

# Row 5 - Chicago
$ws.Range("A5").Value = "12/29/2022 10:46:58"
$ws.Range("B5").Value = "Chicago"
$ws.Range("C5").Value = "52"
$ws.Range("D5").Value = "5%"
$ws.Range("E5").Value = "19 mph"
$ws.Range("F5").Value = "19 mph"

# Row 6 - Nashville
$ws.Range("A6").Value = "12/29/2022 10:47:00"
$ws.Range("B6").Value = "Nashville"
$ws.Range("C6").Value = "60"
$ws.Range("D6").Value = "60%"
$ws.Range("E6").Value = "16 mph"
$ws.Range("F6").Value = "19 mph"

# Restore the original (blank/"General") number format now that the text
# values are locked in, so the cell style index is unchanged.
$ws.Range("C2:D6").NumberFormat = ""
